$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.020.77'
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").Value = '1.630.21'
$ws.Range("E3").Value = '  -0.82%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.99%  '

$ws.Range("E6").Value = '  -1.06%  '

$ws.Range("E7").Value = '  +0.45%  '

$ws.Range("E8").Value = '  -2.78%  '

$ws.Range("E9").Value = '  -3.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.26'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.08%  '

$ws.Range("D12").Value = '1.856.85'
$ws.Range("E12").Value = '  -0.81%  '

$ws.Range("D13").Value = '1.627.38'
$ws.Range("E13").Value = '  -2.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.54%  '

$ws.Range("E15").Value = '  -3.96%  '

$ws.Range("D16").Value = '25.985.60'
$ws.Range("E16").Value = '  -0.25%  '

$ws.Range("E17").Value = '  -3.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.37%  '

$ws.Range("E19").Value = '  +0.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '189.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.05%  '

$ws.Range("E21").Value = '  -3.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.81%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.81%  '

$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("E26").Value = '  -1.24%  '

$ws.Range("E27").Value = '  +0.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.67%  '

$ws.Range("E29").Value = '  -2.93%  '

$ws.Range("E30").Value = '  -1.59%  '

$ws.Range("E31").Value = '  -3.59%  '

$ws.Range("E32").Value = '  -4.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.31%  '

$ws.Range("E34").Value = '  -2.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.47'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.88%  '

$ws.Range("D36").Value = '1.131.08'
$ws.Range("E36").Value = '  -0.37%  '

$ws.Range("E37").Value = '  -6.32%  '

$ws.Range("E38").Value = '  -1.06%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.514'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.99%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0154'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.01'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.41%  '

$ws.Range("E42").Value = '  -3.01%  '

$ws.Range("D43").Value = '1.767.98'
$ws.Range("E43").Value = '  -0.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.33%  '

$ws.Range("E45").Value = '  -2.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '54.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0526'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.21%  '

$ws.Range("E48").Value = '  +0.36%  '

$ws.Range("E49").Value = '  +0.28%  '

$ws.Range("E50").Value = '  +0.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.35%  '
